$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Recipe removed two steps (freezer defrost / scale-into-shape) and
#    replaced the end of the "Scale mix..." step with "...and press
#    into puck shape."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "onto a tray and then into freezer.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "onto a tray and press into puck shape.", 2) | Out-Null

# Remove the now-obsolete "When required, defrost ball..." paragraph.
$d.Content.Find.Execute(
    "When required, defrost ball so that the dough is COMPLETELY at room temperature.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Remove the now-obsolete "Scale into desired shape..." paragraph.
$d.Content.Find.Execute(
    "Scale into desired shape depending on whether for icecream S/W or wagonwheel.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# The two deletions above leave two empty list paragraphs behind; drop them.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "`r") {
        $para.Range.Delete()
    }
}

# ------------------------------------------------------------------
# 2) Choc Chips quantity corrected from 1200g to 1500g.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Choc Chips | 1200g", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Choc Chips | 1500g", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Cool on a cooling wire, not on the tray.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Leave to cool on tray.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Leave to cool on cooling wire.", 2) | Out-Null
